# Apply "Generate Report for Handback" updates:
# - Refresh "Latest HO Xliff Generate Date" (Overview sheet) timestamp
# - Flip zh-cn/de-de "ht" status to "mt"
# - Refresh Correspond Handoff/Handback datetimes for zh-cn and de-de

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date: 2016-08-24 00:15:44 -> 2016-08-24 00:16:28
$wsOverview.Range("G2").Value = "2016-08-24 00:16:28"
$wsOverview.Range("G5").Value = "2016-08-24 00:16:28"
$wsDeDe.Range("H2").Value = "2016-08-24 00:16:28"
$wsDeDe.Range("H5").Value = "2016-08-24 00:16:28"

# Priority: "ht" -> "mt"
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# zh-cn Correspond Handoff Datetime: 2016-08-24 00:15:39 -> 2016-08-24 00:16:23
$wsZhCn.Range("H2").Value = "2016-08-24 00:16:23"
$wsZhCn.Range("H5").Value = "2016-08-24 00:16:23"

# zh-cn Correspond Handback DateTime: 2016-08-24 00:15:55 -> 2016-08-24 00:16:40
$wsZhCn.Range("K2").Value = "2016-08-24 00:16:40"
$wsZhCn.Range("K5").Value = "2016-08-24 00:16:40"

# de-de Correspond Handoff Datetime: 2016-08-24 00:16:06 -> 2016-08-24 00:16:46
$wsDeDe.Range("K2").Value = "2016-08-24 00:16:46"
$wsDeDe.Range("K5").Value = "2016-08-24 00:16:46"
